$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 200; this shifts the existing rows 200-287
# down to 201-288 (matching the target dimension A1:R288).
$ws.Rows("200:200").Insert()

# Populate the newly inserted row 200. Most values mirror the row that
# used to be at 200 (now at 201) except the date (D), min/avg price
# (K/M) and $/Kg (P) columns, which carry the new observation.
$ws.Range("A200").Value2 = 9
$ws.Range("B200").Value2 = "Vega Central Mapocho de Santiago"
$ws.Range("C200").Value2 = "Metropolitana"
$ws.Range("D200").Value2 = 44510
$ws.Range("E200").Value2 = 13
$ws.Range("F200").Value2 = 100112039
$ws.Range("G200").Value2 = "Ciboulette"
$ws.Range("H200").Value2 = "Sin especificar"
$ws.Range("I200").Value2 = "Primera"
$ws.Range("J200").Value2 = 160
$ws.Range("K200").Value2 = 800
$ws.Range("L200").Value2 = 1000
$ws.Range("M200").Value2 = 900
$ws.Range("N200").Value2 = "`$/docena de atados"
$ws.Range("O200").Value2 = "Región Metropolitana"
$ws.Range("P200").Value2 = 300
$ws.Range("Q200").Value2 = 3
$ws.Range("R200").Value2 = "Hortaliza"

# Make sure the D200 cell keeps the date number-format used by the rest
# of the column (style index 2 in the original file).
$ws.Range("D200").NumberFormat = $ws.Range("D201").NumberFormat
